$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6616
$ws.Range("I3").Value = 6919
$ws.Range("D4").Value = 1939
$ws.Range("I4").Value = 1584
$ws.Range("I5").Value = 645
$ws.Range("I6").Value = 7958
$ws.Range("D7").Value = 28129
$ws.Range("I7").Value = 23722

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I6").Value = 98
$ws.Range("I7").Value = 282

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 64
$ws.Range("I7").Value = 270

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I2").Value = 46
$ws.Range("I7").Value = 132

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 244
$ws.Range("I5").Value = 20
$ws.Range("I7").Value = 745

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I6").Value = 277
$ws.Range("I7").Value = 912

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 162
$ws.Range("I6").Value = 162
$ws.Range("I7").Value = 547

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I6").Value = 170
$ws.Range("I7").Value = 748
$ws.Range("I8").Value = 1419
$ws.Range("I11").Value = 361
$ws.Range("I14").Value = 132
$ws.Range("I19").Value = 670
$ws.Range("I20").Value = 588
$ws.Range("I23").Value = 233
$ws.Range("I29").Value = 1432
$ws.Range("I33").Value = 1062
$ws.Range("I34").Value = 106
$ws.Range("I37").Value = 745
$ws.Range("I42").Value = 868
$ws.Range("I47").Value = 171
$ws.Range("I50").Value = 120
$ws.Range("I52").Value = 525
$ws.Range("I54").Value = 479
$ws.Range("I55").Value = 273
$ws.Range("I60").Value = 134
$ws.Range("D63").Value = 326
$ws.Range("I63").Value = 75
$ws.Range("I65").Value = 547
$ws.Range("I67").Value = 912
$ws.Range("I69").Value = 50
$ws.Range("I71").Value = 68
$ws.Range("I76").Value = 343
$ws.Range("I79").Value = 674
$ws.Range("I83").Value = 512
$ws.Range("I85").Value = 1063
$ws.Range("I89").Value = 282
$ws.Range("I90").Value = 309
$ws.Range("I91").Value = 251
$ws.Range("I92").Value = 70
$ws.Range("I94").Value = 242
$ws.Range("I95").Value = 360
$ws.Range("I96").Value = 270
$ws.Range("I97").Value = 193
$ws.Range("D101").Value = 28129
$ws.Range("I101").Value = 23722

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 172
$ws.Range("I7").Value = 512

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I5").Value = 20
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 360

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I6").Value = 342
$ws.Range("I7").Value = 1062

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 103
$ws.Range("I6").Value = 232
$ws.Range("I7").Value = 479

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 497
$ws.Range("I4").Value = 75
$ws.Range("I7").Value = 1432

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I3").Value = 195
$ws.Range("I7").Value = 670

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 156
$ws.Range("I7").Value = 343

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 301
$ws.Range("I6").Value = 276
$ws.Range("I7").Value = 1063

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 263
$ws.Range("I6").Value = 317
$ws.Range("I7").Value = 868

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 83
$ws.Range("I3").Value = 87
$ws.Range("I6").Value = 82
$ws.Range("I7").Value = 273

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 83
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 50

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 91
$ws.Range("I7").Value = 251

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 195
$ws.Range("I3").Value = 219
$ws.Range("I7").Value = 674

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 166
$ws.Range("I3").Value = 167
$ws.Range("I6").Value = 202
$ws.Range("I7").Value = 588

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I4").Value = 38
$ws.Range("I6").Value = 156
$ws.Range("I7").Value = 525

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I4").Value = 11
$ws.Range("I7").Value = 106

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 139
$ws.Range("I7").Value = 242

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 171

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I4").Value = 24
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 143
$ws.Range("I3").Value = 79
$ws.Range("I7").Value = 361

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("I6").Value = 125
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 70

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 422
$ws.Range("I5").Value = 43
$ws.Range("I6").Value = 455
$ws.Range("I7").Value = 1419

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 78
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 309

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I5").Value = 32
$ws.Range("I7").Value = 748
